$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows 2-5 (columns A, B, E, F, G, H, Q, R) are cyclically rotated:
# new row2 = old row5, new row3 = old row2, new row4 = old row3, new row5 = old row4.
# Capture old values first, then write the rotated values back.

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

$oldValues = @{}
foreach ($r in 2..5) {
    $oldValues[$r] = @{}
    foreach ($c in $cols) {
        $oldValues[$r][$c] = $ws.Range("$c$r").Value2
    }
}

# mapping: destination row -> source row
$rowMap = @{ 2 = 5; 3 = 2; 4 = 3; 5 = 4 }

foreach ($destRow in 2..5) {
    $srcRow = $rowMap[$destRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $oldValues[$srcRow][$c]
    }
}
